$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows to append (orderItems), columns A-D.
# D values use the existing shared strings: "paid" / "booked"
$rows = @(
    @(3, 1,   0, "booked"),
    @(4, 45,  0, "booked"),
    @(4, 33,  6, "booked"),
    @(4, 104, 0, "booked"),
    @(5, 33,  6, "paid"),
    @(5, 31,  6, "paid"),
    @(5, 45,  0, "paid"),
    @(5, 28,  0, "paid"),
    @(5, 40,  0, "paid")
)

$startRow = 7
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
